# Commit: "add sheet names to output, indent work done on tab"
#
# 1. Give the three generic "SheetN" tabs real names.
# 2. Switch the active/selected tab from the 2nd sheet to the 3rd sheet,
#    and move the in-sheet selection on that 3rd tab to G31.

$wb = $excel.ActiveWorkbook

# --- 1. Name the sheets -----------------------------------------------
$wb.Worksheets.Item(1).Name = "namedTab1"
$wb.Worksheets.Item(2).Name = "namedTab2"
$wb.Worksheets.Item(3).Name = "namedTab3"

# --- 2. Move the active tab from namedTab2 to namedTab3, with a new ---
#        selection on namedTab3.
$ws3 = $wb.Worksheets.Item("namedTab3")
$ws3.Activate()
$ws3.Range("G31").Select()
